$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.324023666666666
$ws.Range("H2").Value = 3.972071
$ws.Range("I2").Value = 0.01518042398701374
$ws.Range("J2").Value = 0.01518042398701374
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 10.76594625520466
$ws.Range("R2").Value = 96.89351629684198
$ws.Range("S2").Value = 0.0003174219935334893
$ws.Range("T2").Value = 0.0003174219935334893

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.324023666666666
$ws.Range("H3").Value = 3.972071
$ws.Range("I3").Value = 0.01518042398701374
$ws.Range("J3").Value = 0.01518042398701374
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 322.2361117427685
$ws.Range("R3").Value = 2900.125005684916
$ws.Range("S3").Value = 0.009500774623356625
$ws.Range("T3").Value = 0.009500774623356627

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.324023666666666
$ws.Range("H4").Value = 3.972071
$ws.Range("I4").Value = 0.01518042398701374
$ws.Range("J4").Value = 0.01518042398701374
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 137.6225131291072
$ws.Range("R4").Value = 1238.602618161965
$ws.Range("S4").Value = 0.004057647273820572
$ws.Range("T4").Value = 0.004057647273820573

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.324023666666666
$ws.Range("H5").Value = 3.972071
$ws.Range("I5").Value = 0.01518042398701374
$ws.Range("J5").Value = 0.01518042398701374
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 44.24721502773421
$ws.Range("R5").Value = 398.2249352496079
$ws.Range("S5").Value = 0.00130458009630305
$ws.Range("T5").Value = 0.00130458009630305

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 81.17653533333333
$ws.Range("H6").Value = 243.529606
$ws.Range("I6").Value = 0.9307191821270077
$ws.Range("J6").Value = 0.9307191821270075
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 660.0654041046013
$ws.Range("R6").Value = 5940.588636941411
$ws.Range("S6").Value = 0.0194612968954848
$ws.Range("T6").Value = 0.0194612968954848

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 81.17653533333333
$ws.Range("H7").Value = 243.529606
$ws.Range("I7").Value = 0.9307191821270077
$ws.Range("J7").Value = 0.9307191821270075
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 19756.45282566409
$ws.Range("R7").Value = 177808.0754309768
$ws.Range("S7").Value = 0.5824971156660689
$ws.Range("T7").Value = 0.5824971156660689

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.17653533333333
$ws.Range("H8").Value = 243.529606
$ws.Range("I8").Value = 0.9307191821270077
$ws.Range("J8").Value = 0.9307191821270075
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 8437.703253305721
$ws.Range("R8").Value = 75939.32927975149
$ws.Range("S8").Value = 0.2487763289932376
$ws.Range("T8").Value = 0.2487763289932376

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.17653533333333
$ws.Range("H9").Value = 243.529606
$ws.Range("I9").Value = 0.9307191821270077
$ws.Range("J9").Value = 0.9307191821270075
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 2712.818286053143
$ws.Range("R9").Value = 24415.36457447829
$ws.Range("S9").Value = 0.07998444057221632
$ws.Range("T9").Value = 0.07998444057221632

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.192675
$ws.Range("H10").Value = 3.578025
$ws.Range("I10").Value = 0.0136744626508778
$ws.Range("J10").Value = 0.0136744626508778
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 9.697919510949999
$ws.Range("R10").Value = 87.28127559855
$ws.Range("S10").Value = 0.0002859324086635569
$ws.Range("T10").Value = 0.000285932408663557

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.192675
$ws.Range("H11").Value = 3.578025
$ws.Range("I11").Value = 0.0136744626508778
$ws.Range("J11").Value = 0.0136744626508778
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 290.2689462797667
$ws.Range("R11").Value = 2612.420516517901
$ws.Range("S11").Value = 0.00855825817860144
$ws.Range("T11").Value = 0.00855825817860144

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.192675
$ws.Range("H12").Value = 3.578025
$ws.Range("I12").Value = 0.0136744626508778
$ws.Range("J12").Value = 0.0136744626508778
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 123.9697861742084
$ws.Range("R12").Value = 1115.728075567875
$ws.Range("S12").Value = 0.003655111750749636
$ws.Range("T12").Value = 0.003655111750749636

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.192675
$ws.Range("H13").Value = 3.578025
$ws.Range("I13").Value = 0.0136744626508778
$ws.Range("J13").Value = 0.0136744626508778
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 39.85770686113334
$ws.Range("R13").Value = 358.7193617502
$ws.Range("S13").Value = 0.001175160312863169
$ws.Range("T13").Value = 0.001175160312863169

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.525915333333334
$ws.Range("H14").Value = 10.577746
$ws.Range("I14").Value = 0.04042593123510095
$ws.Range("J14").Value = 0.04042593123510094
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 28.67004263952133
$ws.Range("R14").Value = 258.030383755692
$ws.Range("S14").Value = 0.0008453044324763814
$ws.Range("T14").Value = 0.0008453044324763814

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.525915333333334
$ws.Range("H15").Value = 10.577746
$ws.Range("I15").Value = 0.04042593123510095
$ws.Range("J15").Value = 0.04042593123510094
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 858.1245758302464
$ws.Range("R15").Value = 7723.121182472218
$ws.Range("S15").Value = 0.02530085206662018
$ws.Range("T15").Value = 0.02530085206662018

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.525915333333334
$ws.Range("H16").Value = 10.577746
$ws.Range("I16").Value = 0.04042593123510095
$ws.Range("J16").Value = 0.04042593123510094
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 366.4929422866211
$ws.Range("R16").Value = 3298.43648057959
$ws.Range("S16").Value = 0.01080563822249564
$ws.Range("T16").Value = 0.01080563822249564

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.525915333333334
$ws.Range("H17").Value = 10.577746
$ws.Range("I17").Value = 0.04042593123510095
$ws.Range("J17").Value = 0.04042593123510094
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 117.8316806952231
$ws.Range("R17").Value = 1060.485126257008
$ws.Range("S17").Value = 0.003474136513508747
$ws.Range("T17").Value = 0.003474136513508747
